$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = "ALB"
$ws.Range("A7").Value = "CRM"
$ws.Range("A8").Value = "TYL"
$ws.Range("A9").Value = "TMO"
$ws.Range("A11").Value = "F"
$ws.Range("A12").Value = "INCY"
$ws.Range("A13").Value = "SYK"
$ws.Range("A15").Value = "TPR"
$ws.Range("A16").Value = "NXPI"
$ws.Range("A18").Value = "SNPS"
$ws.Range("A14").Value = "SWKS"

$ws.Range("A14").Select()
